# Update the SECTOR_holdings workbook:
#  - bump the "Model holdings provided as of ..." date in the confidentiality
#    footnote from 2021-03-30 to 2021-03-31
#  - refresh the Weight (col D) / Percent Change (col E) figures for rows 2-6

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; unprotect to edit the locked cells, then restore
# protection with the same password once done.
$ws.Unprotect("D382")

$ws.Range("D2").Value = 0.2542595958835814
$ws.Range("E2").Value = -0.00492486425053662

$ws.Range("D3").Value = 0.2511746789489762
$ws.Range("E3").Value = -0.007577965607694637

$ws.Range("D4").Value = 0.2488768529834765
$ws.Range("E4").Value = 0.01591065554960625

$ws.Range("D5").Value = 0.2456888721839659
$ws.Range("E5").Value = 0.003970970833903964

$ws.Range("E6").Value = 0.001779830155025364

$footnote = $ws.Range("A9").Value()
$ws.Range("A9").Value = $footnote.Replace("2021-03-30", "2021-03-31")

$ws.Protect("D382")
